# Stock_Database.xlsx edit:
#  - Insert a new column H ("Current Holdings") that computes a running
#    Buy/Sell balance per stock symbol with SUMIFS, shifting the existing
#    "Price per Unit" column from H to I.
#  - Give the new column a header and fill it down for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Price per Unit" column (H),
# which shifts it to I and leaves a blank H to fill in.
$ws.Columns.Item(8).Insert()

# New header
$ws.Range("H1").Value = "Current Holdings"

# Running Buy/Sell balance per row, based on the stock symbol in column C
for ($r = 2; $r -le 17; $r++) {
    $formula = '=SUMIFS($G$2:G' + $r + ', $C$2:C' + $r + ', C' + $r + ', $F$2:F' + $r + ', "Buy") - SUMIFS($G$2:G' + $r + ', $C$2:C' + $r + ', C' + $r + ', $F$2:F' + $r + ', "Sell")'
    $ws.Range("H$r").Formula = $formula
}

# Re-apply the auto-fit column widths that Excel recalculates after the
# insert (column C narrows slightly, new H + shifted I need sizing).
$ws.Columns.Item(3).ColumnWidth = 20.6
$ws.Columns.Item(8).ColumnWidth = 16
[void]$ws.Range("I17").Select()
